$wb = $excel.ActiveWorkbook

# --- "Videos" sheet: add Times Crossing Cory (R) values and Comments (AC) notes ---
$videos = $wb.Worksheets.Item("Videos")
$videos.Range("R26").Value = 0
$videos.Range("AC26").Value = "4 hours video"

$videos.Range("R27").Value = 1
$videos.Range("AC27").Value = "polyps not out"

$videos.Range("R28").Value = 0
$videos.Range("AC28").Value = "really good deterrence video"

$videos.Range("R29").Value = 1

$videos.Range("R30").Value = 0

$videos.Range("R31").Value = 1
$videos.Range("AC31").Value = "urchin tried crossing"

# --- "Kelp consumption" sheet: insert new column X = IF(V<0.05, 0, V) ---
$ws = $wb.Worksheets.Item("Kelp consumption")
$ws.Columns.Item(24).Insert()

$ws.Cells.Item(1,24).Value = "No zero"
$ws.Range("X2:X77").Formula = "=IF(V2<0.05, 0,V2)"
